$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row updates pulled from the commit diff: (row, new Price text, new Volume(1h) text).
# Price is "" for rows whose D column is unchanged (only the Volume column moved).
$updates = @(
    ,@(2, '58.452.83', '  -3.85%  ')
    ,@(3, '2.981.56', '  -1.34%  ')
    ,@(4, '1.00', '  +0.03%  ')
    ,@(5, '562.49', '  -2.78%  ')
    ,@(6, '134.67', '  +5.72%  ')
    ,@(7, '', '  +0.12%  ')
    ,@(8, '0.518', '  +3.46%  ')
    ,@(9, '2.980.25', '  -1.30%  ')
    ,@(10, '', '  -2.35%  ')
    ,@(11, '', '  -4.96%  ')
    ,@(12, '', '  +1.87%  ')
    ,@(13, '', '  +0.88%  ')
    ,@(14, '33.31', '  +1.58%  ')
    ,@(15, '', '  +0.93%  ')
    ,@(16, '3.471.63', '  -1.19%  ')
    ,@(17, '6.86', '  +6.79%  ')
    ,@(18, '2.974.86', '  -1.34%  ')
    ,@(19, '58.276.04', '  -3.97%  ')
    ,@(20, '425.71', '  -2.19%  ')
    ,@(21, '13.35', '  +1.46%  ')
    ,@(22, '0.693', '  +3.62%  ')
    ,@(23, '', '  -0.04%  ')
    ,@(24, '13.20', '  +2.43%  ')
    ,@(25, '80.10', '  +0.55%  ')
    ,@(26, '', '  -0.04%  ')
    ,@(27, '0.999', '  +0.00%  ')
    ,@(28, '', '  -1.93%  ')
    ,@(29, '', '  +4.35%  ')
    ,@(30, '', '  +5.44%  ')
    ,@(31, '25.58', '  +0.33%  ')
    ,@(32, '6.16', '  -0.85%  ')
    ,@(33, '0.0998', '  +6.20%  ')
    ,@(34, '', '  +1.64%  ')
    ,@(35, '', '  -0.79%  ')
    ,@(36, '0.953', '  -0.74%  ')
    ,@(37, '0.0₃0703', '  +4.92%  ')
    ,@(38, '48.74', '  -3.95%  ')
    ,@(39, '8.78', '  +3.20%  ')
    ,@(40, '', '  +4.13%  ')
    ,@(41, '0.0354', '  -1.88%  ')
    ,@(42, '', '  -1.13%  ')
    ,@(43, '381.60', '  -2.14%  ')
    ,@(44, '2.721.43', '  +2.05%  ')
    ,@(46, '0.244', '  +3.09%  ')
    ,@(47, '123.01', '  +3.15%  ')
    ,@(48, '', '  +2.98%  ')
    ,@(49, '2.02', '  -0.64%  ')
    ,@(50, '23.73', '  -0.21%  ')
    ,@(51, '', '  +1.10%  ')
)

foreach ($u in $updates) {
    $row = $u[0]
    $price = $u[1]
    $volume = $u[2]

    if ($price -ne "") {
        $priceCell = $ws.Range("D$row")
        # A bare numeric-looking string (e.g. "562.49") would be auto-converted to a
        # number by Excel's normal cell-entry parsing, the same as typing it into the
        # grid -- but these Price cells must stay text, exactly like the original
        # inline strings (many prices use "." as a thousands separator, e.g. "58.452.83").
        $looksNumeric = $price -match '^[+-]?[0-9]*\.?[0-9]+$'
        if ($looksNumeric) {
            # Leading apostrophe forces text entry, the same trick used in Excel's UI.
            $priceCell.Value = "'" + $price
            # That apostrophe entry also flips on the quotePrefix cell format; restore
            # the original (unstyled) formatting by copying it from the untouched Link
            # cell in the same row, so only the cell VALUE changes -- matching the diff.
            $priceCell.Style = $ws.Range("C$row").Style
        } else {
            $priceCell.Value = $price
        }
    }

    $ws.Range("E$row").Value = $volume
}
